$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $CellRef, $Val)
    $rng = $Sheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "30.271.02"
Set-TextValue $ws "E2" "  -0.77%  "
Set-TextValue $ws "D3" "2.073.89"
Set-TextValue $ws "E3" "  +2.96%  "
Set-TextValue $ws "D4" "0.9992"
Set-TextValue $ws "D5" "327.25"
Set-TextValue $ws "E5" "  +0.67%  "
Set-TextValue $ws "D6" "0.9989"
Set-TextValue $ws "E6" "  -0.15%  "
Set-TextValue $ws "D7" "0.5191"
Set-TextValue $ws "E7" "  +1.60%  "
Set-TextValue $ws "D8" "0.4315"
Set-TextValue $ws "E8" "  +3.89%  "
Set-TextValue $ws "D9" "0.08826"
Set-TextValue $ws "E9" "  +0.31%  "
Set-TextValue $ws "D10" "46.03"
Set-TextValue $ws "E10" "  +7.39%  "
Set-TextValue $ws "D11" "1.154"
Set-TextValue $ws "E11" "  +1.59%  "
Set-TextValue $ws "D12" "24.24"
Set-TextValue $ws "E12" "  -1.81%  "
Set-TextValue $ws "D13" "2.070.93"
Set-TextValue $ws "E13" "  +2.73%  "
Set-TextValue $ws "D14" "6.650"
Set-TextValue $ws "E14" "  +0.47%  "
Set-TextValue $ws "D15" "7.662"
Set-TextValue $ws "E15" "  +1.61%  "
Set-TextValue $ws "D16" "94.97"
Set-TextValue $ws "E16" "  +0.45%  "
Set-TextValue $ws "D17" "1.000"
Set-TextValue $ws "E17" "  -0.10%  "
Set-TextValue $ws "D18" "0.00001117"
Set-TextValue $ws "E18" "  -0.06%  "
Set-TextValue $ws "D19" "0.06611"
Set-TextValue $ws "E19" "  +1.01%  "
Set-TextValue $ws "D20" "18.74"
Set-TextValue $ws "E20" "  -1.16%  "
Set-TextValue $ws "D21" "0.9987"
Set-TextValue $ws "E21" "  -0.16%  "
Set-TextValue $ws "D22" "6.203"
Set-TextValue $ws "E22" "  -0.50%  "
Set-TextValue $ws "D23" "30.308.78"
Set-TextValue $ws "E23" "  -0.78%  "
Set-TextValue $ws "D24" "12.23"
Set-TextValue $ws "E24" "  +2.61%  "
Set-TextValue $ws "D25" "2.279"
Set-TextValue $ws "E25" "  +2.22%  "
Set-TextValue $ws "D26" "2.315.00"
Set-TextValue $ws "E26" "  +2.89%  "
Set-TextValue $ws "D27" "22.14"
Set-TextValue $ws "E27" "  -0.98%  "
Set-TextValue $ws "D28" "2.544"
Set-TextValue $ws "E28" "  +4.38%  "
Set-TextValue $ws "D29" "161.70"
Set-TextValue $ws "E29" "  -0.97%  "
Set-TextValue $ws "D30" "130.86"
Set-TextValue $ws "E30" "  -0.53%  "
Set-TextValue $ws "D31" "1.183"
Set-TextValue $ws "E31" "  +3.55%  "
Set-TextValue $ws "E32" "  +1.14%  "
Set-TextValue $ws "D33" "1.630"
Set-TextValue $ws "E33" "  +20.62%  "
Set-TextValue $ws "D34" "6.104"
Set-TextValue $ws "E34" "  -0.10%  "
Set-TextValue $ws "D35" "3.820"
Set-TextValue $ws "E35" "  -0.35%  "
Set-TextValue $ws "D36" "0.02576"
Set-TextValue $ws "E36" "  +2.09%  "
Set-TextValue $ws "D37" "9.734"
Set-TextValue $ws "E37" "  +6.25%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws "D38" "12.75"
Set-TextValue $ws "E38" "  +3.37%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D39" "0.06616"
Set-TextValue $ws "E39" "  -0.92%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D40" "5.410"
Set-TextValue $ws "E40" "  -0.96%  "
Set-TextValue $ws "D41" "0.2239"
Set-TextValue $ws "E41" "  +1.80%  "
Set-TextValue $ws "D42" "0.6797"
Set-TextValue $ws "E42" "  +2.10%  "
Set-TextValue $ws "D43" "1.248"
Set-TextValue $ws "E43" "  +1.08%  "
Set-TextValue $ws "D44" "0.9980"
Set-TextValue $ws "E44" "  -0.21%  "
Set-TextValue $ws "D45" "13.91"
Set-TextValue $ws "E45" "  +1.99%  "
Set-TextValue $ws "D46" "0.6339"
Set-TextValue $ws "E46" "  +2.45%  "
Set-TextValue $ws "D47" "2.197"
Set-TextValue $ws "E47" "  -0.61%  "
Set-TextValue $ws "D48" "3.602"
Set-TextValue $ws "E48" "  -1.83%  "
Set-TextValue $ws "D49" "1.232"
Set-TextValue $ws "E49" "  -3.38%  "
Set-TextValue $ws "D50" "1.182"
Set-TextValue $ws "E50" "  +6.67%  "
Set-TextValue $ws "D51" "81.24"
Set-TextValue $ws "E51" "  -0.49%  "
